$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-06-04 Wednesday"; new = "2025-06-05 Thursday"},
    @{old = "521÷7="; new = "792÷9="},
    @{old = "498÷8="; new = "675÷6="},
    @{old = "268÷4="; new = "100÷3="},
    @{old = "597÷3="; new = "865÷9="},
    @{old = "923÷3="; new = "192÷4="},
    @{old = "319÷5="; new = "988÷9="},
    @{old = "119÷5="; new = "486÷3="},
    @{old = "312÷7="; new = "357÷2="},
    @{old = "624÷8="; new = "987÷9="},
    @{old = "742÷6="; new = "373÷6="},
    @{old = "273÷9="; new = "448÷9="},
    @{old = "923÷4="; new = "425÷2="},
    @{old = "554÷9="; new = "805÷3="},
    @{old = "361÷2="; new = "318÷6="},
    @{old = "457÷9="; new = "787÷3="},
    @{old = "376÷8="; new = "932÷5="},
    @{old = "774÷9="; new = "315÷9="},
    @{old = "249÷4="; new = "402÷5="},
    @{old = "138÷9="; new = "295÷3="},
    @{old = "262÷2="; new = "373÷7="},
    @{old = "557÷6="; new = "537÷8="},
    @{old = "844÷6="; new = "603÷5="},
    @{old = "811÷9="; new = "502÷2="},
    @{old = "777÷7="; new = "426÷5="},
    @{old = "801÷6="; new = "501÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
